# cryptos.xlsx price/volume refresh (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price-column values (column D) look like plain decimal
# numbers (e.g. "211.89"), unlike their "28.307.51"-style, thousands-dotted
# neighbours which Excel cannot parse as numbers and so leaves as text.
# Left alone, Excel.Range.Value would auto-coerce those plain-looking values to
# real numbers, but the workbook stores every Price/Volume cell as text.
# Temporarily format them as Text ("@") while writing so the values land as
# strings, then restore the "Normal" cell style so no stray number format
# lingers on the cell (matching the rest of the sheet, which carries none).
$textForceRows = @(5,9,12,16,20,25,26,27,31,34,41,44,45,47,48,51)
foreach ($r in $textForceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '28.307.51'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.573.71'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '211.89'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -4.85%  '
$ws.Range("D9").Value = '23.75'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").Value = '0.0895'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '1.799.66'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").Value = '1.568.67'
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").Value = '0.516'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '28.328.74'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").Value = '7.41'
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '0.0₃0685'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  -1.46%  '
$ws.Range("D25").Value = '2.05'
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("D26").Value = '151.49'
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = '14.94'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").Value = '0.0479'
$ws.Range("E31").Value = '  +2.97%  '
$ws.Range("E32").Value = '  -3.39%  '
$ws.Range("E33").Value = '  -0.61%  '
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").Value = '  -2.00%  '
$ws.Range("D35").Value = '1.383.06'
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("E36").Value = '  +5.45%  '
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").Value = '0.519'
$ws.Range("E41").Value = '  -3.01%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").Value = '0.785'
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").Value = '0.0464'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("E46").Value = '  -4.39%  '
$ws.Range("D47").Value = '0.924'
$ws.Range("E47").Value = '  -5.57%  '
$ws.Range("D48").Value = '62.25'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("D49").Value = '1.711.25'
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").Value = '85.45'
$ws.Range("E51").Value = '  -0.44%  '

# Drop the temporary Text number format again, restoring the default style so
# the cells carry no explicit style index, same as before the edit.
foreach ($r in $textForceRows) {
    $ws.Range("D$r").Style = "Normal"
}
